# Auto-generated script applying the MeteoCat daily-summary refresh
# described by commit: Update automàtic: dades i banners [2026-02-05 12:15]
#
# For each station row, the extraction timestamp (col E) moves forward by
# ~15 min and a handful of the day's running aggregates (cols H, I, J, K, L,
# M, O) are refreshed with newer observed readings. All values are plain
# text in the source sheet (inline/shared strings), so percentage-looking
# values (column H) must be written with a Text number format first --
# otherwise Excel's automatic-type-detection would store "94%" as the
# number 0.94 with a percentage format instead of the literal text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = "2026-02-05 12:14:08"
$ws.Range('H2').NumberFormat = '@'
$ws.Range('H2').Value = "94%"
$ws.Range('I2').Value = "0.3 mm"
$ws.Range('K2').Value = "0.8 MJ/m2"
$ws.Range('M2').Value = "1.6 °C 10:51 TU"
$ws.Range('O2').Value = "-2.0 °C"

$ws.Range('E3').Value = "2026-02-05 12:14:10"

$ws.Range('E4').Value = "2026-02-05 12:14:13"

$ws.Range('E5').Value = "2026-02-05 12:14:16"

$ws.Range('E6').Value = "2026-02-05 12:14:19"
$ws.Range('I6').Value = "0.8 mm"
$ws.Range('K6').Value = "0.7 MJ/m2"

$ws.Range('E7').Value = "2026-02-05 12:14:21"
$ws.Range('H7').NumberFormat = '@'
$ws.Range('H7').Value = "74%"
$ws.Range('I7').Value = "1.1 mm"
$ws.Range('J7').Value = "994.0 hPa"
$ws.Range('K7').Value = "0.7 MJ/m2"
$ws.Range('O7').Value = "9.6 °C"

$ws.Range('E8').Value = "2026-02-05 12:14:24"
$ws.Range('H8').NumberFormat = '@'
$ws.Range('H8').Value = "94%"
$ws.Range('K8').Value = "1.1 MJ/m2"
$ws.Range('M8').Value = "11.9 °C 10:40 TU"
$ws.Range('O8').Value = "5.1 °C"

$ws.Range('E9').Value = "2026-02-05 12:14:27"

$ws.Range('E10').Value = "2026-02-05 12:14:29"
$ws.Range('M10').Value = "8.6 °C 10:59 TU"
$ws.Range('O10').Value = "3.4 °C"

$ws.Range('E11').Value = "2026-02-05 12:14:32"

$ws.Range('E12').Value = "2026-02-05 12:14:35"
$ws.Range('H12').NumberFormat = '@'
$ws.Range('H12').Value = "94%"
$ws.Range('I12').Value = "1.3 mm"
$ws.Range('K12').Value = "0.7 MJ/m2"
$ws.Range('M12').Value = "9.3 °C 10:59 TU"
$ws.Range('O12').Value = "7.3 °C"

$ws.Range('E13').Value = "2026-02-05 12:14:37"
$ws.Range('O13').Value = "4.5 °C"

$ws.Range('E14').Value = "2026-02-05 12:14:40"

$ws.Range('E15').Value = "2026-02-05 12:14:43"
$ws.Range('J15').Value = "993.6 hPa"
$ws.Range('K15').Value = "1.1 MJ/m2"
$ws.Range('L15').Value = "17.6 km/h - 171º 10:53 TU"
$ws.Range('M15').Value = "10.2 °C 10:53 TU"
$ws.Range('O15').Value = "3.3 °C"

$ws.Range('E16').Value = "2026-02-05 12:14:45"
$ws.Range('I16').Value = "0.8 mm"
$ws.Range('K16').Value = "0.6 MJ/m2"
$ws.Range('L16').Value = "22.7 km/h - 241º 10:38 TU"
$ws.Range('O16').Value = "2.5 °C"

$ws.Range('E17').Value = "2026-02-05 12:14:48"

$ws.Range('E18').Value = "2026-02-05 12:14:51"

$ws.Range('E19').Value = "2026-02-05 12:14:54"
$ws.Range('I19').Value = "5.4 mm"
$ws.Range('J19').Value = "994.8 hPa"
$ws.Range('K19').Value = "0.7 MJ/m2"
$ws.Range('M19').Value = "8.4 °C 10:59 TU"
$ws.Range('O19').Value = "5.2 °C"

$ws.Range('E20').Value = "2026-02-05 12:14:56"

$ws.Range('E21').Value = "2026-02-05 12:14:59"

$ws.Range('E22').Value = "2026-02-05 12:15:02"

$ws.Range('E23').Value = "2026-02-05 12:15:05"
$ws.Range('H23').NumberFormat = '@'
$ws.Range('H23').Value = "95%"
$ws.Range('J23').Value = "993.2 hPa"
$ws.Range('K23').Value = "1.2 MJ/m2"
$ws.Range('M23').Value = "7.8 °C 10:44 TU"
$ws.Range('O23').Value = "5.3 °C"

$ws.Range('E24').Value = "2026-02-05 12:15:07"
$ws.Range('H24').NumberFormat = '@'
$ws.Range('H24').Value = "83%"
$ws.Range('J24').Value = "992.3 hPa"
$ws.Range('K24').Value = "1.1 MJ/m2"
$ws.Range('L24').Value = "83.5 km/h - 218º 10:54 TU"
$ws.Range('M24').Value = "10.7 °C 10:55 TU"
$ws.Range('O24').Value = "8.5 °C"

$ws.Range('E25').Value = "2026-02-05 12:15:10"

$ws.Range('E26').Value = "2026-02-05 12:15:13"
$ws.Range('H26').NumberFormat = '@'
$ws.Range('H26').Value = "68%"
$ws.Range('I26').Value = "0.8 mm"
$ws.Range('K26').Value = "1.5 MJ/m2"
$ws.Range('O26').Value = "-2.3 °C"

$ws.Range('E27').Value = "2026-02-05 12:15:15"
$ws.Range('H27').NumberFormat = '@'
$ws.Range('H27').Value = "98%"
$ws.Range('J27').Value = "993.2 hPa"
$ws.Range('K27').Value = "1.2 MJ/m2"
$ws.Range('M27').Value = "11.8 °C 10:53 TU"
$ws.Range('O27').Value = "5.1 °C"

$ws.Range('E28').Value = "2026-02-05 12:15:18"

$ws.Range('E29').Value = "2026-02-05 12:15:21"

$ws.Range('E30').Value = "2026-02-05 12:15:23"

$ws.Range('E31').Value = "2026-02-05 12:15:26"

$ws.Range('E32').Value = "2026-02-05 12:15:29"
$ws.Range('J32').Value = "993.9 hPa"
$ws.Range('K32').Value = "1.0 MJ/m2"
$ws.Range('M32').Value = "13.3 °C 10:56 TU"
$ws.Range('O32').Value = "9.4 °C"

$ws.Range('E33').Value = "2026-02-05 12:15:32"
$ws.Range('H33').NumberFormat = '@'
$ws.Range('H33').Value = "97%"
$ws.Range('M33').Value = "11.9 °C 10:59 TU"
$ws.Range('O33').Value = "4.8 °C"

$ws.Range('E34').Value = "2026-02-05 12:15:34"

$ws.Range('E35').Value = "2026-02-05 12:15:36"
$ws.Range('I35').Value = "2.0 mm"
$ws.Range('K35').Value = "0.9 MJ/m2"

$ws.Range('E36').Value = "2026-02-05 12:15:39"
$ws.Range('H36').NumberFormat = '@'
$ws.Range('H36').Value = "99%"
$ws.Range('J36').Value = "994.7 hPa"
$ws.Range('K36').Value = "2.3 MJ/m2"
$ws.Range('L36').Value = "27.4 km/h - 189º 10:30 TU"
$ws.Range('M36').Value = "14.1 °C 10:59 TU"
$ws.Range('O36').Value = "6.2 °C"
